$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Flow of Events" cell - the paragraph that used to be split across
# many small <w:r> runs ("If system will ", "not ", "found ", "any ", ...)
# collapses down to a single run once the (unchanged) text is re-applied via
# Find & Replace.
# ---------------------------------------------------------------------------
$p1 = "If system will not found any similar data base on proximity hourly, date and position, it update new information and it will assign an ID to the new help request"
$d.Content.Find.Execute($p1, $true, $false, $false, $false, $false, $true, 1, $false, $p1, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Exceptions" cell - "In case of uncertainty" + bookmark + " in
# data, system will insert ..." collapses into one run, and the old _GoBack
# bookmark (which used to sit between the two original runs) goes away - it
# re-appears later, at the very end of the document's last edit (see below).
# ---------------------------------------------------------------------------
$full = "In case of uncertainty in data, system will insert it on DB. This in order to avoid " + [char]0x2018 + "failure to rescue" + [char]0x2019 + " situation "
$d.Content.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, $full, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: "Reached Goals" cell - " [G1]" turns into " [G5] [G6] [G7]". The
# "1" is replaced with "5" in place, then " [G6] [G7]" is appended after the
# closing bracket.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$cell = $t.Cell(7, 2)
$base = $cell.Range.Start

# " [G1]" -> " [G5]"  (nbsp, '[', 'G', '1', ']'  ->  nbsp, '[', 'G', '5', ']')
$d.Range($base + 3, $base + 4).Text = "5"

# Append " [G6] [G7]" right after the "]" and before the paragraph mark.
$insertPoint = $d.Range($base + 5, $base + 5)
$insertPoint.InsertAfter(" [G6] [G7]")

# ---------------------------------------------------------------------------
# Change 4: move the "_GoBack" bookmark from the "Exceptions" cell (removed
# above) to the very end of the "Reached Goals" cell, right after "[G7]".
# Word COM can't collapse-bookmark exactly at end-of-paragraph directly in
# this host, so we temporarily append a throw-away character, anchor the
# bookmark just before it, then remove the throw-away character again - the
# bookmark stays put.
# ---------------------------------------------------------------------------
$cell2 = $t.Cell(7, 2)
$r2 = $cell2.Range
$dummyPoint = $d.Range($r2.End - 1, $r2.End - 1)
$dummyPoint.InsertAfter("Z")

$cell3 = $t.Cell(7, 2)
$r3 = $cell3.Range
$bmPoint = $d.Range($r3.End - 2, $r3.End - 2)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

$cell4 = $t.Cell(7, 2)
$r4 = $cell4.Range
$zRange = $d.Range($r4.End - 2, $r4.End - 1)
$zRange.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Change 5: force the "]" that now sits between " [G5" and " [G6] [G7]" onto
# its own run (matching the three-run split seen in the saved document),
# by toggling a character property on/off which splits the run without
# altering its final formatting.
# ---------------------------------------------------------------------------
$bracketRange = $d.Range($base + 4, $base + 5)
$bracketRange.Font.Bold = 1
$bracketRange.Font.Bold = 0
